$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 4994.2
$ws.Range("I19").Value = 4994.5
$ws.Range("J19").Value = 4994
$ws.Range("K19").Value = 4994.5
$ws.Range("L19").Value = 4994
$ws.Range("M19").Value = -4819.5
$ws.Range("N19").Value = -5344
# Row 32
$ws.Range("H32").Value = 6878.2
$ws.Range("I32").Value = 6697
$ws.Range("K32").Value = 6697
$ws.Range("M32").Value = -6371
# Row 81
$ws.Range("H81").Value = 10328
$ws.Range("J81").Value = 10328
$ws.Range("L81").Value = 10328
$ws.Range("N81").Value = -12324
# Row 84
$ws.Range("H84").Value = 10328
$ws.Range("J84").Value = 10328
$ws.Range("L84").Value = 30984
$ws.Range("N84").Value = -40968
# Row 96
$ws.Range("H96").Value = 720.25
$ws.Range("I96").Value = 672.8
$ws.Range("J96").Value = 754.1429000000001
$ws.Range("K96").Value = 2018.4
$ws.Range("L96").Value = 2262.4287
$ws.Range("M96").Value = -645.3999999999999
$ws.Range("N96").Value = -5008.4287
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
# Row 128
$ws.Range("H128").Value = 68780
$ws.Range("J128").Value = 68780
$ws.Range("L128").Value = 68780
$ws.Range("N128").Value = -78740
# Row 137
$ws.Range("H137").Value = 1142
$ws.Range("I137").Value = 1091.75
$ws.Range("K137").Value = 3275.25
$ws.Range("M137").Value = -725.25

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 1000
$ws.Range("J6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("N6").Value = -1346
# Row 32
$ws.Range("H32").Value = 4634.1577
$ws.Range("I32").Value = 4002.8823
$ws.Range("K32").Value = 4002.8823
$ws.Range("M32").Value = -3715.8823
# Row 74
$ws.Range("H74").Value = 1092.5
$ws.Range("I74").Value = 1092.5
$ws.Range("K74").Value = 1092.5
$ws.Range("M74").Value = -218.5
# Row 77
$ws.Range("H77").Value = 1092.5
$ws.Range("I77").Value = 1092.5
$ws.Range("K77").Value = 5462.5
$ws.Range("M77").Value = -1094.5
# Row 110
$ws.Range("H110").Value = 721.6667
$ws.Range("I110").Value = 666.2
$ws.Range("J110").Value = 999
$ws.Range("K110").Value = 666.2
$ws.Range("L110").Value = 999
$ws.Range("M110").Value = 1378.8
$ws.Range("N110").Value = -5089
# Row 132
$ws.Range("H132").Value = 2409.889
$ws.Range("I132").Value = 2398.625
$ws.Range("K132").Value = 7195.875
$ws.Range("M132").Value = -4665.875

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 3166901.8
$ws.Range("I7").Value = 4750301.5
$ws.Range("J7").Value = 102
$ws.Range("K7").Value = 4750301.5
$ws.Range("L7").Value = 102
$ws.Range("M7").Value = -4750188.5
$ws.Range("N7").Value = -328
# Row 86
$ws.Range("H86").Value = 6240.9165
$ws.Range("I86").Value = 6535.5454
$ws.Range("K86").Value = 6535.5454
$ws.Range("M86").Value = -5412.5454
# Row 89
$ws.Range("H89").Value = 6240.9165
$ws.Range("I89").Value = 6535.5454
$ws.Range("K89").Value = 32677.727
$ws.Range("M89").Value = -27061.727
# Row 105
$ws.Range("H105").Value = 3525
$ws.Range("I105").Value = 4250
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 4250
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -2503
$ws.Range("N105").Value = -6294
# Row 107
$ws.Range("H107").Value = 190
$ws.Range("I107").Value = 190
$ws.Range("K107").Value = 190
$ws.Range("M107").Value = 1730

$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 9097.888999999999
$ws.Range("I12").Value = 666.3077
$ws.Range("J12").Value = 31020
$ws.Range("K12").Value = 666.3077
$ws.Range("L12").Value = 31020
$ws.Range("M12").Value = -496.3077
$ws.Range("N12").Value = -31360
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
# Row 23
$ws.Range("H23").Value = 62002
$ws.Range("J23").Value = 58752.5
$ws.Range("L23").Value = 58752.5
$ws.Range("N23").Value = -59232.5
# Row 27
$ws.Range("H27").Value = 62002
$ws.Range("J27").Value = 58752.5
$ws.Range("L27").Value = 58752.5
$ws.Range("N27").Value = -59136.5
# Row 31
$ws.Range("H31").Value = 2708.889
$ws.Range("I31").Value = 1847.3846
$ws.Range("K31").Value = 1847.3846
$ws.Range("M31").Value = -1552.3846
# Row 34
$ws.Range("H34").Value = 2708.889
$ws.Range("I34").Value = 1847.3846
$ws.Range("K34").Value = 1847.3846
$ws.Range("M34").Value = -1645.3846
# Row 58
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
# Row 99
$ws.Range("H99").Value = 4071.4285
$ws.Range("J99").Value = 4250
$ws.Range("L99").Value = 4250
$ws.Range("N99").Value = -7246
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
# Row 126
$ws.Range("H126").Value = 4071.4285
$ws.Range("J126").Value = 4250
$ws.Range("L126").Value = 12750
$ws.Range("N126").Value = -17690
# Row 134
$ws.Range("H134").Value = 1299.5
$ws.Range("I134").Value = 1299.5
$ws.Range("K134").Value = 3898.5
$ws.Range("M134").Value = -1363.5
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 10254.167
$ws.Range("J13").Value = 12300
$ws.Range("L13").Value = 36900
$ws.Range("N13").Value = -37236
# Row 22
$ws.Range("H22").Value = 127498.75
$ws.Range("I22").Value = 168331.67
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 504995.01
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = -504826.01
$ws.Range("N22").Value = -15338
# Row 27
$ws.Range("H27").Value = 127498.75
$ws.Range("I27").Value = 168331.67
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 504995.01
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = -504893.01
$ws.Range("N27").Value = -15204
# Row 36
$ws.Range("H36").Value = 331.2
$ws.Range("I36").Value = 331.2
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 993.5999999999999
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -824.5999999999999
$ws.Range("N36").ClearContents()
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5444.1665
$ws.Range("I80").Value = 1533
$ws.Range("K80").Value = 1533
$ws.Range("M80").Value = -535
# Row 83
$ws.Range("H83").Value = 5444.1665
$ws.Range("I83").Value = 1533
$ws.Range("K83").Value = 7665
$ws.Range("M83").Value = -2673
# Row 97
$ws.Range("H97").Value = 2443.125
$ws.Range("I97").Value = 2404.3845
$ws.Range("J97").Value = 2611
$ws.Range("K97").Value = 2404.3845
$ws.Range("L97").Value = 2611
$ws.Range("M97").Value = -1908.3845
$ws.Range("N97").Value = -3603
# Row 102
$ws.Range("H102").Value = 1216.8
$ws.Range("I102").Value = 1216.8
$ws.Range("K102").Value = 1216.8
$ws.Range("M102").Value = 405.2
# Row 132
$ws.Range("H132").Value = 4596.4
$ws.Range("I132").Value = 4119
$ws.Range("K132").Value = 12357
$ws.Range("M132").Value = -9827

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 5756.6665
$ws.Range("I16").Value = 5598
$ws.Range("J16").Value = 5915.3335
$ws.Range("K16").Value = 5598
$ws.Range("L16").Value = 5915.3335
$ws.Range("M16").Value = -5428
$ws.Range("N16").Value = -6255.3335
# Row 46
$ws.Range("H46").Value = 3567.647
$ws.Range("I46").Value = 3072.3333
$ws.Range("J46").Value = 4124.875
$ws.Range("K46").Value = 3072.3333
$ws.Range("L46").Value = 4124.875
$ws.Range("M46").Value = -2884.3333
$ws.Range("N46").Value = -4500.875
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 1855.4286
$ws.Range("I13").Value = 998
$ws.Range("J13").Value = 2498.5
$ws.Range("K13").Value = 998
$ws.Range("L13").Value = 2498.5
$ws.Range("M13").Value = -858
$ws.Range("N13").Value = -2778.5
# Row 23
$ws.Range("H23").Value = 757.5
$ws.Range("J23").Value = 2400
$ws.Range("L23").Value = 2400
$ws.Range("N23").Value = -2858
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
# Row 94
$ws.Range("H94").Value = 36217.5
$ws.Range("J94").Value = 36217.5
$ws.Range("L94").Value = 36217.5
$ws.Range("N94").Value = -38019.5
# Row 132
$ws.Range("H132").Value = 2198.6667
$ws.Range("I132").Value = 2098.625
$ws.Range("K132").Value = 6295.875
$ws.Range("M132").Value = -3765.875
# Row 136
$ws.Range("H136").Value = 1487.9
$ws.Range("I136").Value = 1431
$ws.Range("K136").Value = 4293
$ws.Range("M136").Value = -1743
